$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of new row number -> old (current) row number.
# The data rows (2-31) get reshuffled into a new order; this table
# captures, for each destination row, which source row's data must
# land there.
$rowMap = @{2=2; 3=3; 4=6; 5=5; 6=4; 7=7; 8=8; 9=9; 10=10; 11=11; 12=12; 13=17; 14=24; 15=21; 16=20; 17=14; 18=30; 19=23; 20=19; 21=29; 22=22; 23=28; 24=26; 25=15; 26=25; 27=27; 28=16; 29=13; 30=18; 31=31}

# First snapshot every source row's full contents (values + formulas,
# columns A..Z) before any writes happen, since several rows swap with
# each other (e.g. 4<->6) and would otherwise clobber one another.
$snapshots = @{}
foreach ($r in 2..31) {
    $rng = $ws.Range("A" + $r + ":Z" + $r)
    $snapshots[$r] = $rng.Formula
}

# Write each snapshot into its new destination row.
foreach ($newRow in 2..31) {
    $oldRow = $rowMap[$newRow]
    $ws.Range("A" + $newRow + ":Z" + $newRow).Formula = $snapshots[$oldRow]
}

# Bump the "Förändrad" (changed) date column C for every data row.
foreach ($r in 2..31) {
    $ws.Range("C" + $r).Value2 = 46063
}

# Rewriting the wrapped-text "Artnamn" column triggers Excel's
# auto-fit-row-height behaviour; restore the original fixed row height
# (15pt, customHeight) for every data row so rows don't grow taller.
foreach ($r in 2..31) {
    $ws.Rows($r).RowHeight = 15
}
